# Fix: correção para que atualize retorno dos dados para a função de formatação
#
# Duplicates the last "entry block" of rows three times (as if the user
# copy/pasted the same block repeatedly) in four sheets, and refreshes the
# derived totals on the "Estoque" sheet.

function Contains-Col($arr, $col) {
    foreach ($x in $arr) {
        if ($x -eq $col) {
            return $true
        }
    }
    return $false
}

# Copies a block of rows [$startRow, $endRow] (all columns 1..$numCols) to
# $copies successive blocks starting right after $endRow. Value is always
# copied; NumberFormat/HorizontalAlignment are copied too unless the column
# is listed in $noFormatCols (used where the source column intentionally has
# no explicit style, so we don't invent one on the copy).
function Copy-RowBlock($ws, $startRow, $endRow, $numCols, $copies, $noFormatCols) {
    $blockHeight = $endRow - $startRow + 1
    $destRow = $endRow + 1
    for ($c = 0; $c -lt $copies; $c++) {
        for ($r = 0; $r -lt $blockHeight; $r++) {
            $srcRow = $startRow + $r
            for ($col = 1; $col -le $numCols; $col++) {
                $srcCell = $ws.Cells.Item($srcRow, $col)
                $dstCell = $ws.Cells.Item($destRow, $col)
                $dstCell.Value = $srcCell.Value()
                if (-not (Contains-Col $noFormatCols $col)) {
                    $dstCell.NumberFormat = $srcCell.NumberFormat()
                    $dstCell.HorizontalAlignment = $srcCell.HorizontalAlignment()
                }
            }
            $destRow = $destRow + 1
        }
    }
}

$wb = $excel.ActiveWorkbook

$none = @()

# Compras: duplicate rows 27-31 (last purchase entries) three times -> rows 32-46
$wsCompras = $wb.Worksheets.Item("Compras")
Copy-RowBlock $wsCompras 27 31 6 3 $none

# Vendas: duplicate rows 10-13 (last sale entries) three times -> rows 14-25
$wsVendas = $wb.Worksheets.Item("Vendas")
Copy-RowBlock $wsVendas 10 13 7 3 $none

# A Receber: duplicate rows 74-77 (last installment entries) three times -> rows 78-89
$wsReceber = $wb.Worksheets.Item("A Receber")
Copy-RowBlock $wsReceber 74 77 7 3 $none

# A Pagar: duplicate rows 70-73 (last installment entries) three times -> rows 74-85
# Column A ("marcos") has no explicit style in the source rows, so skip format copy there.
$wsPagar = $wb.Worksheets.Item("A Pagar")
Copy-RowBlock $wsPagar 70 73 7 3 @(1)

# Estoque: refresh the totals that derive from the duplicated purchase/sale data
$wsEstoque = $wb.Worksheets.Item("Estoque")
$wsEstoque.Cells.Item(2, 3).Value = 5128
$wsEstoque.Cells.Item(3, 3).Value = 48
$wsEstoque.Cells.Item(4, 3).Value = 48
